# Auto-generated edit script: updates leve-profit calculation columns (H:N)
# across all 8 crafting-job sheets, per the scheduled-runner price refresh.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(55, 8).Value = 380.1111
$ws.Cells.Item(55, 9).Value = 207.5
$ws.Cells.Item(55, 10).Value = 518.2
$ws.Cells.Item(55, 11).Value = 207.5
$ws.Cells.Item(55, 12).Value = 518.2
$ws.Cells.Item(55, 13).Value = 6.5
$ws.Cells.Item(55, 14).Value = -946.2

$ws.Cells.Item(64, 8).Value = 5792
$ws.Cells.Item(64, 9).Value = 5792
$ws.Cells.Item(64, 11).Value = 5792
$ws.Cells.Item(64, 13).Value = -5544

$ws.Cells.Item(67, 8).Value = 5792
$ws.Cells.Item(67, 9).Value = 5792
$ws.Cells.Item(67, 11).Value = 5792
$ws.Cells.Item(67, 13).Value = -4934

$ws.Cells.Item(80, 8).Value = 1937.0625
$ws.Cells.Item(80, 9).Value = 605.8570999999999
$ws.Cells.Item(80, 10).Value = 2972.4443
$ws.Cells.Item(80, 11).Value = 1817.5713
$ws.Cells.Item(80, 12).Value = 8917.332900000001
$ws.Cells.Item(80, 13).Value = -819.5712999999998
$ws.Cells.Item(80, 14).Value = -10913.3329

$ws.Cells.Item(83, 8).Value = 1937.0625
$ws.Cells.Item(83, 9).Value = 605.8570999999999
$ws.Cells.Item(83, 10).Value = 2972.4443
$ws.Cells.Item(83, 11).Value = 5452.7139
$ws.Cells.Item(83, 12).Value = 26751.9987
$ws.Cells.Item(83, 13).Value = -460.7138999999997
$ws.Cells.Item(83, 14).Value = -36735.9987

$ws.Cells.Item(94, 8).Value = 1202
$ws.Cells.Item(94, 9).Value = 641.2
$ws.Cells.Item(94, 10).Value = 4006
$ws.Cells.Item(94, 11).Value = 641.2
$ws.Cells.Item(94, 12).Value = 4006
$ws.Cells.Item(94, 13).Value = -190.2
$ws.Cells.Item(94, 14).Value = -4908

$ws.Cells.Item(103, 8).Value = 476.33334
$ws.Cells.Item(103, 9).Value = 249.5
$ws.Cells.Item(103, 11).Value = 748.5
$ws.Cells.Item(103, 13).Value = -162.5

$ws.Cells.Item(116, 8).Value = 9740
$ws.Cells.Item(116, 9).Value = 9740
$ws.Cells.Item(116, 10).Value = 0
$ws.Cells.Item(116, 11).Value = 9740
$ws.Cells.Item(116, 12).Value = 0
$ws.Cells.Item(116, 13).Value = -6298
$ws.Cells.Item(116, 14).ClearContents()

$ws.Cells.Item(137, 8).Value = 13890324
$ws.Cells.Item(137, 9).Value = 22223284
$ws.Cells.Item(137, 11).Value = 66669852
$ws.Cells.Item(137, 13).Value = -66667302

$ws.Cells.Item(138, 8).Value = 3610.2126
$ws.Cells.Item(138, 9).Value = 4474.077
$ws.Cells.Item(138, 10).Value = 3279.9119
$ws.Cells.Item(138, 11).Value = 13422.231
$ws.Cells.Item(138, 12).Value = 9839.735700000001
$ws.Cells.Item(138, 13).Value = -8282.231
$ws.Cells.Item(138, 14).Value = -20119.7357

$ws.Cells.Item(141, 8).Value = 2702.8
$ws.Cells.Item(141, 9).Value = 2730.5
$ws.Cells.Item(141, 11).Value = 8191.5
$ws.Cells.Item(141, 13).Value = -3011.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(7, 8).Value = 0
$ws.Cells.Item(7, 9).Value = 0
$ws.Cells.Item(7, 11).Value = 0
$ws.Cells.Item(7, 13).ClearContents()

$ws.Cells.Item(102, 8).Value = 4049.5
$ws.Cells.Item(102, 9).Value = 3412.0417
$ws.Cells.Item(102, 11).Value = 3412.0417
$ws.Cells.Item(102, 13).Value = -1790.0417

$ws.Cells.Item(103, 8).Value = 35181
$ws.Cells.Item(103, 10).Value = 35181
$ws.Cells.Item(103, 12).Value = 35181
$ws.Cells.Item(103, 14).Value = -37525

$ws.Cells.Item(122, 8).Value = 5981.706
$ws.Cells.Item(122, 9).Value = 6054.0713
$ws.Cells.Item(122, 11).Value = 18162.2139
$ws.Cells.Item(122, 13).Value = -15712.2139

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(54, 8).Value = 5000
$ws.Cells.Item(54, 9).Value = 5000
$ws.Cells.Item(54, 11).Value = 5000
$ws.Cells.Item(54, 13).Value = -4516

$ws.Cells.Item(105, 8).Value = 6399.75
$ws.Cells.Item(105, 9).Value = 6399.75
$ws.Cells.Item(105, 11).Value = 6399.75
$ws.Cells.Item(105, 13).Value = -4652.75

$ws.Cells.Item(108, 8).Value = 69999
$ws.Cells.Item(108, 10).Value = 69999
$ws.Cells.Item(108, 12).Value = 69999
$ws.Cells.Item(108, 14).Value = -77679

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(5, 8).Value = 2609.8
$ws.Cells.Item(5, 9).Value = 2787.25
$ws.Cells.Item(5, 11).Value = 2787.25
$ws.Cells.Item(5, 13).Value = -2675.25

$ws.Cells.Item(99, 8).Value = 0
$ws.Cells.Item(99, 9).Value = 0
$ws.Cells.Item(99, 10).Value = 0
$ws.Cells.Item(99, 11).Value = 0
$ws.Cells.Item(99, 12).Value = 0
$ws.Cells.Item(99, 13).ClearContents()
$ws.Cells.Item(99, 14).ClearContents()

$ws.Cells.Item(122, 8).Value = 0
$ws.Cells.Item(122, 9).Value = 0
$ws.Cells.Item(122, 11).Value = 0
$ws.Cells.Item(122, 13).ClearContents()

$ws.Cells.Item(126, 8).Value = 0
$ws.Cells.Item(126, 9).Value = 0
$ws.Cells.Item(126, 10).Value = 0
$ws.Cells.Item(126, 11).Value = 0
$ws.Cells.Item(126, 12).Value = 0
$ws.Cells.Item(126, 13).ClearContents()
$ws.Cells.Item(126, 14).ClearContents()

$ws.Cells.Item(134, 8).Value = 4422.25
$ws.Cells.Item(134, 9).Value = 4422.25
$ws.Cells.Item(134, 11).Value = 13266.75
$ws.Cells.Item(134, 13).Value = -10731.75

$ws.Cells.Item(141, 8).Value = 396998.1
$ws.Cells.Item(141, 10).Value = 396998.1
$ws.Cells.Item(141, 12).Value = 396998.1
$ws.Cells.Item(141, 14).Value = -407358.1

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(40, 8).Value = 48.142857
$ws.Cells.Item(40, 9).Value = 47.833332
$ws.Cells.Item(40, 10).Value = 50
$ws.Cells.Item(40, 11).Value = 191.333328
$ws.Cells.Item(40, 12).Value = 200
$ws.Cells.Item(40, 13).Value = -122.333328
$ws.Cells.Item(40, 14).Value = -338

$ws.Cells.Item(57, 8).Value = 16833.334
$ws.Cells.Item(57, 9).Value = 1000
$ws.Cells.Item(57, 10).Value = 20000
$ws.Cells.Item(57, 11).Value = 3000
$ws.Cells.Item(57, 12).Value = 60000
$ws.Cells.Item(57, 13).Value = -2441
$ws.Cells.Item(57, 14).Value = -61118

$ws.Cells.Item(119, 8).Value = 4295485.5
$ws.Cells.Item(119, 9).Value = 4295485.5
$ws.Cells.Item(119, 11).Value = 12886456.5
$ws.Cells.Item(119, 13).Value = -12881618.5

$ws.Cells.Item(137, 8).Value = 7443.5557
$ws.Cells.Item(137, 9).Value = 6333
$ws.Cells.Item(137, 11).Value = 18999
$ws.Cells.Item(137, 13).Value = -13899

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(102, 8).Value = 2480.5715
$ws.Cells.Item(102, 9).Value = 2480.5715
$ws.Cells.Item(102, 11).Value = 2480.5715
$ws.Cells.Item(102, 13).Value = -858.5715

$ws.Cells.Item(122, 8).Value = 9411.375
$ws.Cells.Item(122, 9).Value = 9960.25
$ws.Cells.Item(122, 10).Value = 8862.5
$ws.Cells.Item(122, 11).Value = 29880.75
$ws.Cells.Item(122, 12).Value = 26587.5
$ws.Cells.Item(122, 13).Value = -27430.75
$ws.Cells.Item(122, 14).Value = -31487.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(46, 8).Value = 2000
$ws.Cells.Item(46, 9).Value = 0
$ws.Cells.Item(46, 11).Value = 0
$ws.Cells.Item(46, 13).ClearContents()

$ws.Cells.Item(68, 8).Value = 2497
$ws.Cells.Item(68, 9).Value = 2497
$ws.Cells.Item(68, 11).Value = 2497
$ws.Cells.Item(68, 13).Value = -1748

$ws.Cells.Item(71, 8).Value = 2497
$ws.Cells.Item(71, 9).Value = 2497
$ws.Cells.Item(71, 11).Value = 12485
$ws.Cells.Item(71, 13).Value = -8741

$ws.Cells.Item(100, 8).Value = 1659.4
$ws.Cells.Item(100, 9).Value = 1699.25
$ws.Cells.Item(100, 10).Value = 1500
$ws.Cells.Item(100, 11).Value = 1699.25
$ws.Cells.Item(100, 12).Value = 1500
$ws.Cells.Item(100, 13).Value = -1158.25
$ws.Cells.Item(100, 14).Value = -2582

$ws.Cells.Item(122, 8).Value = 1998.8334
$ws.Cells.Item(122, 9).Value = 2397.6
$ws.Cells.Item(122, 11).Value = 7192.799999999999
$ws.Cells.Item(122, 13).Value = -4742.799999999999

$ws.Cells.Item(130, 8).Value = 76712.5
$ws.Cells.Item(130, 10).Value = 76712.5
$ws.Cells.Item(130, 12).Value = 76712.5
$ws.Cells.Item(130, 14).Value = -86752.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(51, 8).Value = 157766.45
$ws.Cells.Item(51, 9).Value = 157766.45
$ws.Cells.Item(51, 11).Value = 157766.45
$ws.Cells.Item(51, 13).Value = -157256.45

$ws.Cells.Item(131, 8).Value = 0
$ws.Cells.Item(131, 10).Value = 0
$ws.Cells.Item(131, 12).Value = 0
$ws.Cells.Item(131, 14).ClearContents()

$ws.Cells.Item(132, 9).Value = 1126.2222
$ws.Cells.Item(132, 11).Value = 3378.6666
$ws.Cells.Item(132, 13).Value = -848.6665999999996
